$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mejoras")

# Add the new improvement entry (row 6, columns C and D)
$ws.Range("C6").Value = "frmPedido"
$ws.Range("D6").Value = "Al buscar un comobo y si no existe aún así deshabilita los botones de buscar y los texBox"

# Update selection to D7 as recorded in the saved view state
$ws.Range("D7").Select()
